$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing 2020/2021 values in row 5 (M5, N5, O5) ---
$ws.Range("M5").Value = 2.6
$ws.Range("N5").Value = 2.4
$ws.Range("O5").Value = 3.3

# --- Add new column P for year 2022 ---

# P3: thin separator row cell just under the title - mirror O3's bottom border
$p3 = $ws.Range("P3")
$p3.Borders.Item(9).LineStyle = 1
$p3.Borders.Item(9).Weight = -4138
$p3.Borders.Item(9).Color = 0

# P4: year header "2022" - mirror O4's number format / alignment / bottom border
$p4 = $ws.Range("P4")
$p4.Value = 2022
$p4.NumberFormat = "0"
$p4.HorizontalAlignment = -4152
$p4.VerticalAlignment = -4108
$p4.Borders.Item(9).LineStyle = 1
$p4.Borders.Item(9).Weight = -4138
$p4.Borders.Item(9).Color = 0

# P5: data value for 2022 - mirror O5's top+bottom border and vertical alignment
$p5 = $ws.Range("P5")
$p5.Value = 2.6
$p5.VerticalAlignment = -4108
$p5.Borders.Item(8).LineStyle = 1
$p5.Borders.Item(8).Weight = -4138
$p5.Borders.Item(8).Color = 0
$p5.Borders.Item(9).LineStyle = 1
$p5.Borders.Item(9).Weight = -4138
$p5.Borders.Item(9).Color = 0

# Move the selection to P3 (matches the workbook's saved cursor position)
$ws.Range("P3").Select()
